# Natmi following Dr Hou advice
#
# Re-computes the Icam2 -> Itgam ligand-receptor report: a new source
# cluster ("sCs") is now included among the sending/target clusters
# (alongside ECs, FAPs, M2), and every expression/specificity statistic is
# refreshed for all sending x target cluster combinations. This grows the
# table from 6 data rows (3 senders x 2 targets) to 8 data rows
# (4 senders x 2 targets), so the sheet dimension grows from A1:T7 to
# A1:T9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ncols = 20
$startRow = 2

# Flat row-major array: 8 data rows x 20 columns (A:T).
# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression
# value, Ligand total expression value, Ligand derived specificity of
# average expression value, Ligand derived specificity of total expression
# value, Receptor-expressing cells, Receptor detection rate, Receptor
# average expression value, Receptor total expression value, Receptor
# derived specificity of average expression value, Receptor derived
# specificity of total expression value, Edge average expression weight,
# Edge total expression weight, Edge average expression derived
# specificity, Edge total expression derived specificity.
$data = @(
  "ECs", "Icam2", "Itgam", "M2", 3, 1, 32.270062, 96.810186, 0.8981134838283896, 0.8981134838283895, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 1482.206415061592, 13339.85773555433, 0.8868167443172602, 0.88681674431726,
  "ECs", "Icam2", "Itgam", "sCs", 3, 1, 32.270062, 96.810186, 0.8981134838283896, 0.8981134838283895, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 18.88112722270134, 169.930145004312, 0.01129673951112956, 0.01129673951112956,
  "FAPs", "Icam2", "Itgam", "M2", 3, 1, 1.529744666666667, 4.589234, 0.04257457924772188, 0.04257457924772188, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 70.26318568398135, 632.3686711558321, 0.04203906348026308, 0.04203906348026308,
  "FAPs", "Icam2", "Itgam", "sCs", 3, 1, 1.529744666666667, 4.589234, 0.04257457924772188, 0.04257457924772188, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 0.8950495251475556, 8.055445726328001, 0.0005355157674587997, 0.0005355157674587997,
  "M2", "Icam2", "Itgam", "M2", 3, 1, 1.822405666666667, 5.467217, 0.05071967640590832, 0.05071967640590832, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 83.70549055585734, 753.349415002716, 0.0500817091748587, 0.05008170917485869,
  "M2", "Icam2", "Itgam", "sCs", 3, 1, 1.822405666666667, 5.467217, 0.05071967640590832, 0.05071967640590832, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 1.066284695818222, 9.596562262363999, 0.0006379672310496254, 0.0006379672310496254,
  "sCs", "Icam2", "Itgam", "M2", 2, 0.6666666666666666, 0.308728, 0.926184, 0.008592260517980134, 0.008592260517980134, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 14.180283326048, 127.622549934432, 0.008484184500159283, 0.008484184500159281,
  "sCs", "Icam2", "Itgam", "sCs", 2, 0.6666666666666666, 0.308728, 0.926184, 0.008592260517980134, 0.008592260517980134, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 0.1806359295253333, 1.625723365728, 0.0001080760178208522, 0.0001080760178208522
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $startRow + [math]::Floor($i / $ncols)
  $c = ($i % $ncols) + 1
  $ws.Cells.Item($r, $c).Value = $data[$i]
}
